# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) is re-sorted from descending
# (2307..2301) to ascending (2301..2307). The "Valor Mora" of 34666 used
# to be tied to period 2307 (row 16); after the re-sort, period 2307 is
# the last row (22), so the 34666 figure moves down with it, while the
# vacated first row (now period 2301) takes the standard 40000 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for E16:E22
$ws.Range("E16").Value = "2301"
$ws.Range("E17").Value = "2302"
$ws.Range("E18").Value = "2303"
$ws.Range("E19").Value = "2304"
$ws.Range("E20").Value = "2305"
$ws.Range("E21").Value = "2306"
$ws.Range("E22").Value = "2307"

# Swap the "Valor Mora" amounts so the 34666 value moves from period 2307
# (row 16) down to period 2307's new location (row 22); every other row
# keeps the standard 40000 value.
$ws.Range("F16").Value = 40000
$ws.Range("F22").Value = 34666
